$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 92.62231230044235
$ws.Range("C2").Value = 92.13548435394372
$ws.Range("D2").Value = 90.01252588952323
$ws.Range("E2").Value = 92.51996281739414

$ws.Range("B3").Value = 98.49886966425875
$ws.Range("C3").Value = 97.46262661849083
$ws.Range("D3").Value = 98.39605930508836
$ws.Range("E3").Value = 98.1113544923256

$ws.Range("B4").Value = 99.22649049131337
$ws.Range("C4").Value = 99.16043794794793
$ws.Range("D4").Value = 99.24435836521205
$ws.Range("E4").Value = 99.26908377696874

$ws.Range("B5").Value = 98.71270779674394
$ws.Range("C5").Value = 98.72924675671565
$ws.Range("D5").Value = 98.7149701986567
$ws.Range("E5").Value = 98.69679794927404

$ws.Range("B6").Value = 98.30911481604721
$ws.Range("C6").Value = 98.20710235099992
$ws.Range("D6").Value = 98.23620060404036
$ws.Range("E6").Value = 98.18310202875765

$ws.Range("B7").Value = 97.25502922678261
$ws.Range("C7").Value = 97.25175620981535
$ws.Range("D7").Value = 97.29418695927653
$ws.Range("E7").Value = 97.27983962192604

$ws.Range("B8").Value = 95.89108503202644
$ws.Range("C8").Value = 95.86849104052054
$ws.Range("D8").Value = 95.87431958803579
$ws.Range("E8").Value = 95.82466499915455
